$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11.08709081533996
$ws.Range("B2").Value = [double]"1.110223024625157e-16"
$ws.Range("C2").Value = 0.003769768628071014
$ws.Range("D2").Value = 0.5078544707115397
$ws.Range("E2").Value = 0.2579161634216982

$ws.Range("A3").Value = 8.526064092341917
$ws.Range("B3").Value = [double]"1.110223024625157e-16"
$ws.Range("C3").Value = 0.004230512232555903
$ws.Range("D3").Value = 0.5699247785938388
$ws.Range("E3").Value = 0.3248142532552362

$ws.Range("A4").Value = 10.30892457946977
$ws.Range("B4").Value = [double]"1.110223024625157e-16"
$ws.Range("C4").Value = 0.003729885883511057
$ws.Range("D4").Value = 0.5024815600299138
$ws.Range("E4").Value = 0.2524877181700959

$ws.Range("A5").Value = 9.590471825901743
$ws.Range("B5").Value = [double]"1.110223024625157e-16"
$ws.Range("C5").Value = 0.005761314904354511
$ws.Range("D5").Value = 0.776150957797813
$ws.Range("E5").Value = 0.6024103092904624

$ws.Range("A6").Value = 7.81511277691239
$ws.Range("B6").Value = [double]"1.110223024625157e-16"
$ws.Range("C6").Value = 0.003657881049432447
$ws.Range("D6").Value = 0.4927812360823468
$ws.Range("E6").Value = 0.2428333466348455
